$d = $word.ActiveDocument

$d.Content.Find.Execute("85×90=7650", $true, $false, $false, $false, $false, $true, 1, $false, "35×22=770", 2) | Out-Null
$d.Content.Find.Execute("69×85=5865", $true, $false, $false, $false, $false, $true, 1, $false, "43×31=1333", 2) | Out-Null
$d.Content.Find.Execute("96×18=1728", $true, $false, $false, $false, $false, $true, 1, $false, "21×61=1281", 2) | Out-Null
$d.Content.Find.Execute("90×21=1890", $true, $false, $false, $false, $false, $true, 1, $false, "24×37=888", 2) | Out-Null
$d.Content.Find.Execute("71×66=4686", $true, $false, $false, $false, $false, $true, 1, $false, "69×67=4623", 2) | Out-Null
$d.Content.Find.Execute("21×81=1701", $true, $false, $false, $false, $false, $true, 1, $false, "92×13=1196", 2) | Out-Null
$d.Content.Find.Execute("12×75=900", $true, $false, $false, $false, $false, $true, 1, $false, "30×31=930", 2) | Out-Null
$d.Content.Find.Execute("16×55=880", $true, $false, $false, $false, $false, $true, 1, $false, "83×84=6972", 2) | Out-Null
$d.Content.Find.Execute("56×64=3584", $true, $false, $false, $false, $false, $true, 1, $false, "70×75=5250", 2) | Out-Null
$d.Content.Find.Execute("43×82=3526", $true, $false, $false, $false, $false, $true, 1, $false, "43×23=989", 2) | Out-Null
$d.Content.Find.Execute("20×82=1640", $true, $false, $false, $false, $false, $true, 1, $false, "93×54=5022", 2) | Out-Null
$d.Content.Find.Execute("36×40=1440", $true, $false, $false, $false, $false, $true, 1, $false, "47×57=2679", 2) | Out-Null
$d.Content.Find.Execute("29×81=2349", $true, $false, $false, $false, $false, $true, 1, $false, "93×50=4650", 2) | Out-Null
$d.Content.Find.Execute("62×95=5890", $true, $false, $false, $false, $false, $true, 1, $false, "23×61=1403", 2) | Out-Null
$d.Content.Find.Execute("13×64=832", $true, $false, $false, $false, $false, $true, 1, $false, "15×97=1455", 2) | Out-Null
$d.Content.Find.Execute("59×93=5487", $true, $false, $false, $false, $false, $true, 1, $false, "41×45=1845", 2) | Out-Null
$d.Content.Find.Execute("53×98=5194", $true, $false, $false, $false, $false, $true, 1, $false, "31×12=372", 2) | Out-Null
$d.Content.Find.Execute("54×65=3510", $true, $false, $false, $false, $false, $true, 1, $false, "20×82=1640", 2) | Out-Null
$d.Content.Find.Execute("14×81=1134", $true, $false, $false, $false, $false, $true, 1, $false, "76×38=2888", 2) | Out-Null
$d.Content.Find.Execute("85×97=8245", $true, $false, $false, $false, $false, $true, 1, $false, "35×26=910", 2) | Out-Null
$d.Content.Find.Execute("67×79=5293", $true, $false, $false, $false, $false, $true, 1, $false, "81×43=3483", 2) | Out-Null
$d.Content.Find.Execute("28×87=2436", $true, $false, $false, $false, $false, $true, 1, $false, "98×99=9702", 2) | Out-Null
$d.Content.Find.Execute("16×92=1472", $true, $false, $false, $false, $false, $true, 1, $false, "48×67=3216", 2) | Out-Null
$d.Content.Find.Execute("28×77=2156", $true, $false, $false, $false, $false, $true, 1, $false, "78×43=3354", 2) | Out-Null
$d.Content.Find.Execute("90×98=8820", $true, $false, $false, $false, $false, $true, 1, $false, "65×31=2015", 2) | Out-Null
